# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Updates the "Metadata" sheet (Version/Date bump, Publisher + new
# Jurisdiction row replacing the old duplicated Contact row) and the
# "Elements" sheet (Short/Definition text for the root Extension row,
# which now documents the SourceFileId extension instead of a generic
# "An Extension" placeholder).

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ---------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date bump
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# The old sheet had two identical "Contact" / "No display for
# ContactDetail" rows (10 and 11). Drop the duplicate (row 11); the rows
# below shift up by one.
$wsMeta.Rows.Item(11).Delete()

# Publisher now has a value, and what used to be the "Contact" row (10)
# becomes a "Jurisdiction" row.
$wsMeta.Range("B9").Value = "Alvearie Team"
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# ---- Elements sheet ----------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Short / Definition for the root Extension row (row 2) now describe the
# SourceFileId extension specifically instead of the generic placeholder
# text.
$wsElem.Range("K2").Value = "Source File ID"
$wsElem.Range("L2").Value = "The ID for a file from which the data producer or data integrator extracted knowledge, to produce the data within this FHIR resource or element"
